$wb = $excel.ActiveWorkbook

# RTMF-passengers sheet: update the LDVs row (row 2) mode-shift fractions.
$ws = $wb.Worksheets.Item("RTMF-passengers")

# HDVs column (C) and rail column (E) fractions change.
$ws.Range("C2").Value = 0.15
$ws.Range("E2").Value = 0.05

# The "Non-motorized/eliminated" column (I2) was previously a live formula
# (=1-SUM(B2:G2)); it is now a hard-coded value instead.
$ws.Range("I2").Value = 0.8

# Make this the active sheet/selection, matching the saved view state.
$ws.Activate()
$ws.Range("E4").Select()
